$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-12 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-13 Wednesday", 2)

$d.Content.Find.Execute("438÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "464÷2=", 2)
$d.Content.Find.Execute("671÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "363÷7=", 2)
$d.Content.Find.Execute("996÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "849÷3=", 2)
$d.Content.Find.Execute("647÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "506÷8=", 2)
$d.Content.Find.Execute("319÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "533÷8=", 2)

$d.Content.Find.Execute("164÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "203÷5=", 2)
$d.Content.Find.Execute("301÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "285÷7=", 2)
$d.Content.Find.Execute("888÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷3=", 2)
$d.Content.Find.Execute("835÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "296÷8=", 2)
$d.Content.Find.Execute("158÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "116÷3=", 2)

$d.Content.Find.Execute("680÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "644÷6=", 2)
$d.Content.Find.Execute("923÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "810÷3=", 2)
$d.Content.Find.Execute("449÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "111÷5=", 2)
$d.Content.Find.Execute("771÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "767÷3=", 2)
$d.Content.Find.Execute("336÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "393÷2=", 2)

$d.Content.Find.Execute("850÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "381÷9=", 2)
$d.Content.Find.Execute("825÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "495÷9=", 2)
$d.Content.Find.Execute("558÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "819÷2=", 2)
$d.Content.Find.Execute("512÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "634÷4=", 2)
$d.Content.Find.Execute("206÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷4=", 2)

$d.Content.Find.Execute("313÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "258÷6=", 2)
$d.Content.Find.Execute("405÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "653÷2=", 2)
$d.Content.Find.Execute("677÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "659÷4=", 2)
$d.Content.Find.Execute("204÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "929÷2=", 2)
$d.Content.Find.Execute("943÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "115÷7=", 2)
